$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values before overwriting, since the new layout is a
# cyclic rotation of rows 2, 5 and 6. (Note: use Value2, not Value --
# bare .Value resolves to the property's method-group descriptor here.)
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2

# Row 2 becomes what used to be row 6 (org:resource / str)
$ws.Range("A2").Value2 = $a6
$ws.Range("B2").Value2 = $b6

# Row 5 becomes what used to be row 2 (stream:datastream / dict)
$ws.Range("A5").Value2 = $a2
$ws.Range("B5").Value2 = $b2

# Row 6 becomes what used to be row 5 (operation_end_time / datetime)
$ws.Range("A6").Value2 = $a5
$ws.Range("B6").Value2 = $b5
